$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H74").Value = 5000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 5000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H98").Value = 8178.923
$ws.Range("I98").Value = 3213
$ws.Range("J98").Value = 13972.5
$ws.Range("K98").Value = 3213
$ws.Range("L98").Value = 13972.5
$ws.Range("M98").Value = -1715
$ws.Range("N98").Value = -16968.5

$ws.Range("H122").Value = 8178.923
$ws.Range("I122").Value = 3213
$ws.Range("J122").Value = 13972.5
$ws.Range("K122").Value = 9639
$ws.Range("L122").Value = 41917.5
$ws.Range("M122").Value = -7189
$ws.Range("N122").Value = -46817.5

$ws.Range("H125").Value = 1931
$ws.Range("I125").Value = 1396.5
$ws.Range("K125").Value = 12568.5
$ws.Range("M125").Value = -10108.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3534.5557
$ws.Range("I32").Value = 3534.5557
$ws.Range("K32").Value = 3534.5557
$ws.Range("M32").Value = -3247.5557

$ws.Range("H61").Value = 4995
$ws.Range("I61").Value = 4990
$ws.Range("K61").Value = 4990
$ws.Range("M61").Value = -4778

$ws.Range("H74").Value = 2173.8823
$ws.Range("I74").Value = 1425.4286
$ws.Range("K74").Value = 1425.4286
$ws.Range("M74").Value = -551.4286

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

$ws.Range("H77").Value = 2173.8823
$ws.Range("I77").Value = 1425.4286
$ws.Range("K77").Value = 7127.143
$ws.Range("M77").Value = -2759.143

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

$ws.Range("H97").Value = 630
$ws.Range("I97").Value = 575.7143
$ws.Range("K97").Value = 575.7143
$ws.Range("M97").Value = -79.71429999999998

$ws.Range("H132").Value = 619
$ws.Range("I132").Value = 619
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1857
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 673
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 4995
$ws.Range("I136").Value = 4990
$ws.Range("K136").Value = 14970
$ws.Range("M136").Value = -12420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8125
$ws.Range("I86").Value = 4500
$ws.Range("K86").Value = 4500
$ws.Range("M86").Value = -3377

$ws.Range("H89").Value = 8125
$ws.Range("I89").Value = 4500
$ws.Range("K89").Value = 22500
$ws.Range("M89").Value = -16884

$ws.Range("H105").Value = 1899
$ws.Range("I105").Value = 1899
$ws.Range("K105").Value = 1899
$ws.Range("M105").Value = -152

$ws.Range("H107").Value = 1066.6666
$ws.Range("I107").Value = 1066.6666
$ws.Range("K107").Value = 1066.6666
$ws.Range("M107").Value = 853.3334

$ws.Range("H134").Value = 1974.75
$ws.Range("I134").Value = 1974.75
$ws.Range("K134").Value = 5924.25
$ws.Range("M134").Value = -3389.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 41.2
$ws.Range("I7").Value = 42
$ws.Range("J7").Value = 34
$ws.Range("K7").Value = 42
$ws.Range("L7").Value = 34
$ws.Range("M7").Value = 71
$ws.Range("N7").Value = -260

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H99").Value = 1400
$ws.Range("I99").Value = 1400
$ws.Range("K99").Value = 1400
$ws.Range("M99").Value = 98

$ws.Range("H105").Value = 1296.5
$ws.Range("I105").Value = 1294
$ws.Range("K105").Value = 1294
$ws.Range("M105").Value = 453

$ws.Range("H107").Value = 526.3333
$ws.Range("I107").Value = 336
$ws.Range("J107").Value = 716.6667
$ws.Range("K107").Value = 336
$ws.Range("L107").Value = 716.6667
$ws.Range("M107").Value = 1584
$ws.Range("N107").Value = -4556.6667

$ws.Range("H126").Value = 1400
$ws.Range("I126").Value = 1400
$ws.Range("K126").Value = 4200
$ws.Range("M126").Value = -1730

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 37.18182
$ws.Range("J12").Value = 32.285713
$ws.Range("L12").Value = 96.857139
$ws.Range("N12").Value = -442.857139

$ws.Range("H13").Value = 3975
$ws.Range("I13").Value = 4200
$ws.Range("K13").Value = 12600
$ws.Range("M13").Value = -12432

$ws.Range("H75").Value = 1000
$ws.Range("J75").Value = 1000
$ws.Range("L75").Value = 3000
$ws.Range("N75").Value = -4996

$ws.Range("H78").Value = 1000
$ws.Range("J78").Value = 1000
$ws.Range("L78").Value = 9000
$ws.Range("N78").Value = -18984

$ws.Range("H118").Value = 4000
$ws.Range("I118").Value = 4000
$ws.Range("K118").Value = 12000
$ws.Range("M118").Value = -10757

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 712
$ws.Range("I107").Value = 68.5
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 68.5
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = 1851.5
$ws.Range("N107").Value = -5839

$ws.Range("H134").Value = 4500
$ws.Range("J134").Value = 4500
$ws.Range("L134").Value = 13500
$ws.Range("N134").Value = -18570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null

$ws.Range("H68").Value = 3439.8
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 4333
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 4333
$ws.Range("M68").Value = -1351
$ws.Range("N68").Value = -5831

$ws.Range("H71").Value = 3439.8
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 4333
$ws.Range("K71").Value = 10500
$ws.Range("L71").Value = 21665
$ws.Range("M71").Value = -6756
$ws.Range("N71").Value = -29153

$ws.Range("H100").Value = 1665.6666
$ws.Range("I100").Value = 100
$ws.Range("J100").Value = 2448.5
$ws.Range("K100").Value = 100
$ws.Range("L100").Value = 2448.5
$ws.Range("M100").Value = 441
$ws.Range("N100").Value = -3530.5

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1034.5
$ws.Range("I81").Value = 1046
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 2092
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -1031
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 1034.5
$ws.Range("I84").Value = 1046
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 10460
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -5156
